$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the weekly price block (rows 147-148),
# pushing the existing data (previously rows 147-164) down to rows 149-166.
$ws.Range("A147:T148").EntireRow.Insert()

# New row 147: Limón Sutil De Gase, week of 44522
$ws.Range("A147").Value = 1
$ws.Range("B147").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C147").Value = "Arica y Parinacota"
$ws.Range("D147").Value = 44522
$ws.Range("E147").Value = 15
$ws.Range("F147").Value = "Fruta"
$ws.Range("G147").Value = 100102
$ws.Range("H147").Value = "Cítricos"
$ws.Range("I147").Value = 100102003
$ws.Range("J147").Value = "Limón"
$ws.Range("K147").Value = "Sutil De Gase"
$ws.Range("L147").Value = "Primera"
$ws.Range("M147").Value = 200
$ws.Range("N147").Value = 20000
$ws.Range("O147").Value = 21000
$ws.Range("P147").Value = 20500
$ws.Range("Q147").Value = "$/caja 24 kilos"
$ws.Range("R147").Value = "Perú"
$ws.Range("S147").Value = 854
$ws.Range("T147").Value = 24

# New row 148: Limón Tahití, week of 44522
$ws.Range("A148").Value = 1
$ws.Range("B148").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C148").Value = "Arica y Parinacota"
$ws.Range("D148").Value = 44522
$ws.Range("E148").Value = 15
$ws.Range("F148").Value = "Fruta"
$ws.Range("G148").Value = 100102
$ws.Range("H148").Value = "Cítricos"
$ws.Range("I148").Value = 100102003
$ws.Range("J148").Value = "Limón"
$ws.Range("K148").Value = "Tahití"
$ws.Range("L148").Value = "Primera"
$ws.Range("M148").Value = 200
$ws.Range("N148").Value = 24000
$ws.Range("O148").Value = 25000
$ws.Range("P148").Value = 24500
$ws.Range("Q148").Value = "$/caja 24 kilos"
$ws.Range("R148").Value = "Perú"
$ws.Range("S148").Value = 1021
$ws.Range("T148").Value = 24
